# Westeros MACRO input workbook: add a "commodity" column to the
# "config" sheet (inserted before the existing "level" column) and
# update the active sheet / selection to match the re-saved workbook.

$wb = $excel.ActiveWorkbook

# --- "MERtoPPP" sheet: this was the previously-active sheet; its
# selection moves to a single cell (K2) and it stops being the
# tab-selected sheet once we activate "config" below.
$wsMer = $wb.Worksheets.Item("MERtoPPP")
$wsMer.Activate()
$wsMer.Range("K2").Select()

# --- "config" sheet: insert a new column C ("commodity"), shifting
# the existing "level"/"useful" column from C to D.
$ws = $wb.Worksheets.Item("config")
$ws.Activate()

$ws.Columns.Item(3).Insert()

$ws.Cells.Item(1, 3).Value = "commodity"
$ws.Cells.Item(2, 3).Value = "light"

# Match the column width Excel computed when it auto-fit the new
# "commodity" column (displayed width 11 characters). ColumnWidth is
# specified in pre-padding character units, so back out the ~5/6
# character gridline padding Excel adds to the stored <col width>.
$ws.Columns.Item(3).ColumnWidth = 10.166666666666666

# Final selection/cursor position on the "config" sheet.
$ws.Range("A3").Select()
